# fix(FN-3460): fix invalid facility utilisation values in e2e report fixtures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 800000
$ws.Range("F2").Value = 761579.37

# Row 3
$ws.Range("E3").Value = 800000
$ws.Range("F3").Value = 761579.37

# Row 4
$ws.Range("E4").Value = 800000
$ws.Range("F4").Value = 761579.37

# Row 5
$ws.Range("E5").Value = 800000
$ws.Range("F5").Value = 761579.37
$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

# Row 6
$ws.Range("E6").Value = 800000
$ws.Range("G6").Value = 761579.37

# Resize columns G:H to fit the new (wider) values, matching Excel's auto bestFit behaviour
$ws.Range("G:H").EntireColumn.AutoFit()

# Update the visible selection to match the authored state
$ws.Range("E2:H6").Select()
